# Final commit for assignment.
# Fills in the "Preconditions" (E), "Method Inputs" (F) and
# "Expected Result" (G) columns for test cases 1-8 (rows 7-14) of the
# unit test plan, moves the active selection to G14, and widens column F
# to fit the newly-entered text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - __init__ / balance greater than overdraft limit
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "(12345,678,1000.00,date(2024,10,1),300,0.10)"
$ws.Range("G7").Value = "(12345,678,1000.00,2024-10-1,300,0.10)"

# Row 8 - __init__ / overdraft limit has invalid type.
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "(12345,678,1000.00,date(2024,10,1),invalid type,0.10)"
$ws.Range("G8").Value = "(12345,678,1000.00,date(2024,10,1),-100,0.10)"

# Row 9 - __init__ / overdraft rate has invalid type.
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "(12345,678,1000.00,date(2024,10,1),300,invalid rate)"
$ws.Range("G9").Value = "(12345,678,1000.00,date(2024,10,1),300,0.05)"

# Row 10 - __init__ / date created has invalid type
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "(12345,678,1000.00,invalid date,300,0.10)"
$ws.Range("G10").Value = "(12345,678,1000.00,date.today(),300,0.10)"

# Row 11 - get_service_charges / balance less than overdraft limit
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = "(12345,678,1000.00,date(2024,10,1),300,0.10)"
$ws.Range("G11").Value = "service charge = base service charge"

# Row 12 - get_service_charges / balance equal to overdraft limit
$ws.Range("E12").Value = "None"
$ws.Range("F12").Value = "(12345,678,200.00,date(2024,10,1),300,0.10)"
$ws.Range("G12").Value = "service charge is correctly calculated base on the fomula"

# Row 13 - get_service_charges / appropriate value returned based on attribute values.
$ws.Range("E13").Value = "None"
$ws.Range("F13").Value = "(12345,678,300.00,date(2024,10,1),300,0.10)"
$ws.Range("G13").Value = "service charge equal to base service charge"

# Row 14 - __str__ / Attributes are set to input values
$ws.Range("E14").Value = "None"
$ws.Range("F14").Value = "str(self.chequeaccount1)"
$ws.Range("G14").Value = "(`"Account Number: 12345 Balance: `$1000.00`"`n                        `"\nOverdraft Limit: `$300.00 Overdraft Rate: 10.00%`"`n                        `"\nAccount Type: Chequing`")"

# E10 and E12 previously had the "no data yet" border style; now that they
# hold the same "None" input as the rest of column E, align their format
# with the other Method Inputs cells (style copied from E7).
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Widen column F to accommodate the newly-entered text.
$ws.Columns.Item(6).ColumnWidth = 41

# Move the active selection to G14, matching the saved workbook state.
$ws.Range("G14").Select() | Out-Null
